$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J59
$iVals = @(7,4,8,5,6,6,9,6,8,9,7,7,5,3,10,5,7,9,8,8,9,7,7,8,6,7,8,9,7,8,9,9,9,8,11,9,5,5,4,6,10,5,9,8,8,6,8,9,8,9,9,9,7,5,6,9,5,5)
$jVals = @(7,5,8,5,6,7,9,6,8,9,7,7,5,3,10,6,7,10,8,8,9,7,7,8,7,8,8,9,7,8,9,9,9,8,11,9,6,5,5,6,10,6,9,8,8,6,8,9,8,9,9,9,7,5,6,9,5,5)

for ($r = 2; $r -le 59; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
